$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - copy formatting (bold/border/alignment) from the
# neighboring header cell G1 so the new column matches the existing header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Corresponding numeric value for the new column in the data row.
$ws.Range("H2").Value = 0
